$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Near the end of the document: remove the duplicate bold
#    "Play 15 Armadillos Slot Free - Features Exciting Bonus Rounds"
#    paragraph entirely, and replace the text of the following italic
#    paragraph with the new image-generation prompt.
#
#    Do this FIRST, while the "Read our review..." sentence is still
#    unique in the document (before we add a second copy of it below).
# ---------------------------------------------------------------------

$count = $d.Paragraphs.Count
$titleDupe = $d.Paragraphs($count - 1)
if ($titleDupe.Range.Text -notmatch "Play 15 Armadillos Slot Free") {
    throw "Unexpected second-to-last paragraph content: " + $titleDupe.Range.Text
}
$titleDupe.Range.Delete()

$promptText = @'
Prompt: Create a cartoon-style feature image for the game "15 Armadillos" that features a happy Maya warrior with glasses. For the feature image of "15 Armadillos", let's have a cartoon-style design featuring a happy Maya warrior with glasses. The warrior can be depicted wearing a headdress made of colorful feathers, with intricate designs on their face and body. They can be holding a staff or weapon made of stone or wood, with a happy expression on their face. In the background, we can see the Everglades National Park with its lush greenery and animals like alligators and otters. The image can be bright and colorful to reflect the fun and adventurous nature of the game.
'@

$found = $d.Content.Find.Execute(
    "Read our review and play 15 Armadillos slot for free. Enjoy exciting bonus rounds such as Wild Respins, Armadillo Link, and Free Spins.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $promptText, 2)

if (-not $found) {
    throw "Could not find the closing italic paragraph text to replace"
}

# ---------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
#
#    We use Range.InsertXML on the 2nd paragraph (the Heading2 "Gameplay
#    and How to Play the 15 Armadillos Slot" paragraph), replacing it
#    with [new Normal-style "Meta description" paragraph] + [the
#    original Heading2 paragraph unchanged]. This lets us add a brand
#    new paragraph with no <w:pPr> (i.e. Normal style, matching the
#    diff) without leaving any w:rsid* side effects behind, and without
#    disturbing the existing Heading2 paragraph that follows it.
# ---------------------------------------------------------------------

$wordmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$metaParaXml = @"
<w:p $wordmlNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review and play 15 Armadillos slot for free. Enjoy exciting bonus rounds such as Wild Respins, Armadillo Link, and Free Spins.</w:t></w:r></w:p>
"@

$gameplayHeading = $d.Paragraphs(2)
if ($gameplayHeading.Range.Text -notmatch "Gameplay and How to Play") {
    throw "Unexpected paragraph 2 content: " + $gameplayHeading.Range.Text
}

$gameplayXml = @"
<w:p $wordmlNs><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Gameplay and How to Play the 15 Armadillos Slot</w:t></w:r></w:p>
"@

$insertionRange = $gameplayHeading.Range
$insertionRange.Collapse(1)
$insertionRange.InsertXML($metaParaXml + $gameplayXml)

Write-Output "done"
